$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 20-23: copy the formatting (bold/centered/bordered) from the
# last existing labeled row (A19) before writing values, so no duplicate
# style entries get created in styles.xml ---
$ws.Range("A19").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A: label text (every row from 1 onward shifts because four new
# strings were inserted earlier in the shared-string table; rows 20-23
# are brand new rows appended at the end)
$ws.Range("A1").Value  = "___269882__any_any_any_any"
$ws.Range("A2").Value  = "Surveyed habs"
$ws.Range("A3").Value  = "▷ approved"
$ws.Range("A4").Value  = "Surveyed habs cat-III"
$ws.Range("A5").Value  = "DPR habs"
$ws.Range("A6").Value  = "DPR habs cat-III"
$ws.Range("A7").Value  = "STATUS (1 Approved habs)"
$ws.Range("A8").Value  = "▷ Completed"
$ws.Range("A9").Value  = "▷ Ongoing"
$ws.Range("A10").Value = "▷ Completed cat-III"
$ws.Range("A11").Value = "▷ Cert submitted 📑"
$ws.Range("A12").Value = "▷ non-surveyed"
$ws.Range("A13").Value = "Progress Qty"
$ws.Range("A14").Value = "▷ DTR_100"
$ws.Range("A15").Value = "▷ DTR_63"
$ws.Range("A16").Value = "▷ DTR_25"
$ws.Range("A17").Value = "▷ HT"
$ws.Range("A18").Value = "▷ HT_CONDUCTOR"
$ws.Range("A19").Value = "▷ LT_3P"
$ws.Range("A20").Value = "▷ LT_1P"
$ws.Range("A21").Value = "▷ POLE_HT_8M"
$ws.Range("A22").Value = "▷ POLE_LT_8M"
$ws.Range("A23").Value = "▷ POLE_9M"

# Column B: values
$ws.Range("B2").Value  = 1
$ws.Range("B5").Value  = 2
$ws.Range("B6").Value  = 2
$ws.Range("B7").Value  = "▪️▪️▪️"
$ws.Range("B9").Value  = 1
$ws.Range("B10").Value = 1
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 1
$ws.Range("B13").Value = "▪️▪️▪️"
$ws.Range("B15").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("B19").Value = 0.6
$ws.Range("B20").Value = 0.45
$ws.Range("B21").Value = 0
$ws.Range("B22").Value = 26
$ws.Range("B23").Value = 0
